$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = 4 (Price), Column E = 5 (Volume(1h))
# Some new Price values are plain decimal numbers (e.g. "575.58") which Excel
# would otherwise auto-convert to a numeric value instead of keeping them as
# text (matching the source data, which stores prices as text). Force those
# specific cells to text format before assigning so they stay strings.

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "69.725.27"
$ws.Cells.Item(2, 5).Value = "  +0.69%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "2.507.37"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.08%  "

# Row 5 - BNB
Set-TextValue 5 4 "575.58"
$ws.Cells.Item(5, 5).Value = "  +0.14%  "

# Row 6 - Solana
Set-TextValue 6 4 "166.99"
$ws.Cells.Item(6, 5).Value = "  +0.68%  "

# Row 8 - XRP
$ws.Cells.Item(8, 5).Value = "  -0.08%  "

# Row 9 - LidoStakedEther
$ws.Cells.Item(9, 4).Value = "2.508.46"
$ws.Cells.Item(9, 5).Value = "  +0.49%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 5).Value = "  +4.01%  "

# Row 11 - TRON
$ws.Cells.Item(11, 5).Value = "  -0.04%  "

# Row 12 - Cardano
Set-TextValue 12 4 "0.355"
$ws.Cells.Item(12, 5).Value = "  +4.27%  "

# Row 13 - Toncoin
$ws.Cells.Item(13, 5).Value = "  +2.40%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "2.970.49"
$ws.Cells.Item(14, 5).Value = "  +0.56%  "

# Row 15 - ShibaInu
$ws.Cells.Item(15, 5).Value = "  +3.10%  "

# Row 16 - WrappedBTC
$ws.Cells.Item(16, 4).Value = "69.594.62"
$ws.Cells.Item(16, 5).Value = "  +0.59%  "

# Row 17 - Avalanche
Set-TextValue 17 4 "24.86"
$ws.Cells.Item(17, 5).Value = "  +0.74%  "

# Row 18 - WrappedEther
$ws.Cells.Item(18, 4).Value = "2.529.04"
$ws.Cells.Item(18, 5).Value = "  +1.22%  "

# Row 19 - Chainlink
Set-TextValue 19 4 "11.23"
$ws.Cells.Item(19, 5).Value = "  -1.08%  "

# Row 20 - Uniswap
$ws.Cells.Item(20, 5).Value = "  -3.02%  "

# Row 21 - BitcoinCash
Set-TextValue 21 4 "349.28"
$ws.Cells.Item(21, 5).Value = "  +0.73%  "

# Row 22 - Polkadot
$ws.Cells.Item(22, 5).Value = "  -0.23%  "

# Row 23 - SuiNetwork
$ws.Cells.Item(23, 5).Value = "  +0.99%  "

# Row 24 - Dai
$ws.Cells.Item(24, 5).Value = "  +0.01%  "

# Row 25 - Litecoin
Set-TextValue 25 4 "70.48"
$ws.Cells.Item(25, 5).Value = "  +3.57%  "

# Row 26 - NEARProtocol
$ws.Cells.Item(26, 5).Value = "  +0.09%  "

# Row 27 - Aptos
$ws.Cells.Item(27, 5).Value = "  -0.52%  "

# Row 28 - WrappedeETH
$ws.Cells.Item(28, 4).Value = "2.646.11"
$ws.Cells.Item(28, 5).Value = "  +0.31%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue 29 4 "0.997"
$ws.Cells.Item(29, 5).Value = "  +0.07%  "

# Row 30 - PEPE
$ws.Cells.Item(30, 4).Value = "0.0₃0894"
$ws.Cells.Item(30, 5).Value = "  +0.21%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue 31 4 "7.86"
$ws.Cells.Item(31, 5).Value = "  +0.68%  "

# Row 32 - Bittensor
Set-TextValue 32 4 "458.81"
$ws.Cells.Item(32, 5).Value = "  -1.68%  "

# Row 33 - Fetch.AI
$ws.Cells.Item(33, 5).Value = "  -2.60%  "

# Row 34 - PancakeSwap
$ws.Cells.Item(34, 5).Value = "  -0.22%  "

# Row 36 - was Monero, now Kaspa
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 36 4 "0.117"
$ws.Cells.Item(36, 5).Value = "  +2.21%  "

# Row 37 - was Kaspa, now Monero
$ws.Cells.Item(37, 2).Value = "Monero"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 37 4 "158.36"
$ws.Cells.Item(37, 5).Value = "  +3.57%  "

# Row 38 - WhiteBITCoin
$ws.Cells.Item(38, 5).Value = "  +0.71%  "

# Row 39 - EthereumClassic
Set-TextValue 39 4 "18.50"
$ws.Cells.Item(39, 5).Value = "  +0.91%  "

# Row 40 - USDe
$ws.Cells.Item(40, 5).Value = "  +0.06%  "

# Row 41 - PolygonEcosystemToken
$ws.Cells.Item(41, 5).Value = "  +1.61%  "

# Row 42 - RenderToken
$ws.Cells.Item(42, 5).Value = "  -0.49%  "

# Row 43 - Stacks
$ws.Cells.Item(43, 5).Value = "  +1.66%  "

# Row 44 - OKB
Set-TextValue 44 4 "38.05"
$ws.Cells.Item(44, 5).Value = "  +0.00%  "

# Row 45 - dogwifhat
Set-TextValue 45 4 "2.20"
$ws.Cells.Item(45, 5).Value = "  -3.77%  "

# Row 46 - ImmutableX
$ws.Cells.Item(46, 5).Value = "  -5.97%  "

# Row 47 - Aave
Set-TextValue 47 4 "141.53"
$ws.Cells.Item(47, 5).Value = "  -0.55%  "

# Row 48 - Filecoin
Set-TextValue 48 4 "3.48"
$ws.Cells.Item(48, 5).Value = "  -0.34%  "

# Row 49 - ARBITRUM
$ws.Cells.Item(49, 5).Value = "  -0.71%  "

# Row 50 - Cronos
$ws.Cells.Item(50, 5).Value = "  +0.68%  "

# Row 51 - Mantle
$ws.Cells.Item(51, 5).Value = "  -0.50%  "
